$d = $word.ActiveDocument
$t = $d.Tables(1)

# Update the date heading (first paragraph, above the table)
$d.Paragraphs(1).Range.Find.Execute("2023-10-06 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-07 Saturday", 2) | Out-Null

# Row 1
$t.Cell(1,1).Range.Text = "30÷5=6, 0"
$t.Cell(1,2).Range.Text = "69÷7=9, 6"
$t.Cell(1,3).Range.Text = "11÷5=2, 1"
$t.Cell(1,4).Range.Text = "30÷7=4, 2"
$t.Cell(1,5).Range.Text = "35÷7=5, 0"

# Row 2 (table row 5)
$t.Cell(5,1).Range.Text = "42÷2=21, 0"
$t.Cell(5,2).Range.Text = "53÷8=6, 5"
$t.Cell(5,3).Range.Text = "93÷5=18, 3"
$t.Cell(5,4).Range.Text = "55÷9=6, 1"
$t.Cell(5,5).Range.Text = "37÷8=4, 5"

# Row 3 (table row 9)
$t.Cell(9,1).Range.Text = "29÷4=7, 1"
$t.Cell(9,2).Range.Text = "63÷4=15, 3"
$t.Cell(9,3).Range.Text = "68÷7=9, 5"
$t.Cell(9,4).Range.Text = "11÷5=2, 1"
$t.Cell(9,5).Range.Text = "79÷4=19, 3"

# Row 4 (table row 13)
$t.Cell(13,1).Range.Text = "83÷6=13, 5"
$t.Cell(13,2).Range.Text = "77÷2=38, 1"
$t.Cell(13,3).Range.Text = "55÷8=6, 7"
$t.Cell(13,4).Range.Text = "75÷9=8, 3"
$t.Cell(13,5).Range.Text = "53÷5=10, 3"

# Row 5 (table row 17) - one answer removed, one inserted; net cell count unchanged,
# so every cell's text is simply reassigned to its final value.
$t.Cell(17,1).Range.Text = "40÷9=4, 4"
$t.Cell(17,2).Range.Text = "56÷9=6, 2"
$t.Cell(17,3).Range.Text = "51÷4=12, 3"
$t.Cell(17,4).Range.Text = "65÷5=13, 0"
$t.Cell(17,5).Range.Text = "47÷9=5, 2"
